$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 29.01.2022 12:45"

# Update D3: was text "+0.6" -> numeric 0.6
$ws.Range("D3").Value = 0.6

# Update E3: was text "2022-01-29 12:30:12" -> numeric date serial 44590.52097222222
$ws.Range("E3").Value = 44590.52097222222
$ws.Range("E3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
